$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "30.370.75"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.847.20"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "233.64"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.4674"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").Value = "0.2728"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "0.06291"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").Value = "1.839.37"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").Value = "0.07465"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "16.25"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "4.940"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "83.93"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").Value = "0.6201"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "30.322.47"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "229.97"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D19").Value = "0.000007325"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "12.37"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "4.920"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "5.871"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "9.170"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "165.42"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Value = "17.85"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "1.871"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "0.1028"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "4.083"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").Value = "3.812"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").Value = "0.04888"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "1.143"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "0.7059"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").Value = "2.698"
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "0.01892"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "2.664"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "0.8717"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "105.59"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").Value = "0.9994"
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "5.525"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").Value = "0.4024"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "7.081"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "61.46"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "8.653"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("D48").Value = "33.34"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").Value = "0.05516"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").Value = "1.344"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").Value = "0.3644"
$ws.Range("E51").Value = "  -1.56%  "
